$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 50 data rows to 102; new rows (51:102) need the same
# bold/bordered/centered style already used by the existing column-A label
# cells (copy format from A2, which carries that style).
$ws.Range("A2").Copy()
$ws.Range("A51:A102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-sorted (ascending by count/frequency in column B) label/value pairs.
$ws.Cells.Item(2, 1).Value = 'Best Cinematography'
$ws.Cells.Item(2, 2).Value = 0.0
$ws.Cells.Item(3, 1).Value = ' Best Sound Mixing'
$ws.Cells.Item(3, 2).Value = 0.0
$ws.Cells.Item(4, 1).Value = 'sports'
$ws.Cells.Item(4, 2).Value = 0.0
$ws.Cells.Item(5, 1).Value = ' Best Supporting Actor'
$ws.Cells.Item(5, 2).Value = 0.0
$ws.Cells.Item(6, 1).Value = ' Best Visual Effects'
$ws.Cells.Item(6, 2).Value = 0.0
$ws.Cells.Item(7, 1).Value = ' Cinematography'
$ws.Cells.Item(7, 2).Value = 0.0
$ws.Cells.Item(8, 1).Value = 'documentary'
$ws.Cells.Item(8, 2).Value = 0.0
$ws.Cells.Item(9, 1).Value = 'Writing (adapted screenplay)'
$ws.Cells.Item(9, 2).Value = 0.0
$ws.Cells.Item(10, 1).Value = 'western'
$ws.Cells.Item(10, 2).Value = 0.0
$ws.Cells.Item(11, 1).Value = 'Visual Effects'
$ws.Cells.Item(11, 2).Value = 0.0
$ws.Cells.Item(12, 1).Value = ' Sound Editing'
$ws.Cells.Item(12, 2).Value = 0.0
$ws.Cells.Item(13, 1).Value = ' Sound Mixing'
$ws.Cells.Item(13, 2).Value = 0.0
$ws.Cells.Item(14, 1).Value = 'Costume Design'
$ws.Cells.Item(14, 2).Value = 0.0
$ws.Cells.Item(15, 1).Value = 'Best Supporting Actress (Patricia Arquette)'
$ws.Cells.Item(15, 2).Value = 0.0
$ws.Cells.Item(16, 1).Value = 'Best Supporting Actor (Jared Leto)'
$ws.Cells.Item(16, 2).Value = 0.0
$ws.Cells.Item(17, 1).Value = 'Best Supporting Actor'
$ws.Cells.Item(17, 2).Value = 0.0
$ws.Cells.Item(18, 1).Value = 'Best Film Editing'
$ws.Cells.Item(18, 2).Value = 0.0
$ws.Cells.Item(19, 1).Value = ' Production Design'
$ws.Cells.Item(19, 2).Value = 0.0
$ws.Cells.Item(20, 1).Value = ' Best Film Editing (Tom Cross)'
$ws.Cells.Item(20, 2).Value = 0.0
$ws.Cells.Item(21, 1).Value = 'Best Actress (Cate Blanchett)'
$ws.Cells.Item(21, 2).Value = 0.0
$ws.Cells.Item(22, 1).Value = 'fantasy'
$ws.Cells.Item(22, 2).Value = 0.0
$ws.Cells.Item(23, 1).Value = 'sci-fi'
$ws.Cells.Item(23, 2).Value = 0.0
$ws.Cells.Item(24, 1).Value = '_may'
$ws.Cells.Item(24, 2).Value = 0.0
$ws.Cells.Item(25, 1).Value = '_march'
$ws.Cells.Item(25, 2).Value = 0.0
$ws.Cells.Item(26, 1).Value = 'period'
$ws.Cells.Item(26, 2).Value = 0.0
$ws.Cells.Item(27, 1).Value = 'crime'
$ws.Cells.Item(27, 2).Value = 0.0
$ws.Cells.Item(28, 1).Value = ' Best Actor (Matthew McConaughey)'
$ws.Cells.Item(28, 2).Value = 0.0
$ws.Cells.Item(29, 1).Value = 'mystery'
$ws.Cells.Item(29, 2).Value = 0.0
$ws.Cells.Item(30, 1).Value = '_april'
$ws.Cells.Item(30, 2).Value = 0.0
$ws.Cells.Item(31, 1).Value = 'remake'
$ws.Cells.Item(31, 2).Value = 0.0
$ws.Cells.Item(32, 1).Value = '_january'
$ws.Cells.Item(32, 2).Value = 0.0
$ws.Cells.Item(33, 1).Value = 'family'
$ws.Cells.Item(33, 2).Value = 0.0
$ws.Cells.Item(34, 1).Value = ' Best Cinematography'
$ws.Cells.Item(34, 2).Value = 0.0
$ws.Cells.Item(35, 1).Value = '_february'
$ws.Cells.Item(35, 2).Value = 0.0
$ws.Cells.Item(36, 1).Value = '_august'
$ws.Cells.Item(36, 2).Value = 0.00006357761134831232
$ws.Cells.Item(37, 1).Value = 'Best Supporting Actor (J.K. Simmons)'
$ws.Cells.Item(37, 2).Value = 0.0001140432457568003
$ws.Cells.Item(38, 1).Value = 'thrilled'
$ws.Cells.Item(38, 2).Value = 0.0001668234954020425
$ws.Cells.Item(39, 1).Value = ' Best Makeup and Hairstyling (Adruitha Lee and Robin Mathews)'
$ws.Cells.Item(39, 2).Value = 0.000380144152522669
$ws.Cells.Item(40, 1).Value = 'horror'
$ws.Cells.Item(40, 2).Value = 0.0003824117451385328
$ws.Cells.Item(41, 1).Value = ' Best Original Screenplay'
$ws.Cells.Item(41, 2).Value = 0.0004145747755971396
$ws.Cells.Item(42, 1).Value = 'adventure'
$ws.Cells.Item(42, 2).Value = 0.0004172408513244529
$ws.Cells.Item(43, 1).Value = 'musical'
$ws.Cells.Item(43, 2).Value = 0.0007123563873162332
$ws.Cells.Item(44, 1).Value = '_october'
$ws.Cells.Item(44, 2).Value = 0.00114411205725164
$ws.Cells.Item(45, 1).Value = 'action'
$ws.Cells.Item(45, 2).Value = 0.001186355060138209
$ws.Cells.Item(46, 1).Value = 'Actor in a Leading Role'
$ws.Cells.Item(46, 2).Value = 0.001420115005115009
$ws.Cells.Item(47, 1).Value = '_july'
$ws.Cells.Item(47, 2).Value = 0.001717328583726326
$ws.Cells.Item(48, 1).Value = '_september'
$ws.Cells.Item(48, 2).Value = 0.001746772754914211
$ws.Cells.Item(49, 1).Value = 'war'
$ws.Cells.Item(49, 2).Value = 0.001779618008305992
$ws.Cells.Item(50, 1).Value = ' Directing'
$ws.Cells.Item(50, 2).Value = 0.00192229757209605
$ws.Cells.Item(51, 1).Value = 'Best Actress'
$ws.Cells.Item(51, 2).Value = 0.002020073346014945
$ws.Cells.Item(52, 1).Value = ' Best Sound'
$ws.Cells.Item(52, 2).Value = 0.002358188125460815
$ws.Cells.Item(53, 1).Value = '_june'
$ws.Cells.Item(53, 2).Value = 0.002437418473197094
$ws.Cells.Item(54, 1).Value = ' Writing (original screenplay)'
$ws.Cells.Item(54, 2).Value = 0.002871346710184615
$ws.Cells.Item(55, 1).Value = ' Best Animated Feature'
$ws.Cells.Item(55, 2).Value = 0.00288908928879292
$ws.Cells.Item(56, 1).Value = 'Best Editing'
$ws.Cells.Item(56, 2).Value = 0.002895040039130526
$ws.Cells.Item(57, 1).Value = 'animation'
$ws.Cells.Item(57, 2).Value = 0.002988747715731908
$ws.Cells.Item(58, 1).Value = 'sequel'
$ws.Cells.Item(58, 2).Value = 0.003071738319954514
$ws.Cells.Item(59, 1).Value = ' Best Adapted Screenplay'
$ws.Cells.Item(59, 2).Value = 0.003521041232027668
$ws.Cells.Item(60, 1).Value = 'romantic'
$ws.Cells.Item(60, 2).Value = 0.00381649901292658
$ws.Cells.Item(61, 1).Value = 'biography'
$ws.Cells.Item(61, 2).Value = 0.003862568355657896
$ws.Cells.Item(62, 1).Value = 'Animated Feature Film'
$ws.Cells.Item(62, 2).Value = 0.004581733264141152
$ws.Cells.Item(63, 1).Value = 'original'
$ws.Cells.Item(63, 2).Value = 0.004629564027753592
$ws.Cells.Item(64, 1).Value = 'comedy'
$ws.Cells.Item(64, 2).Value = 0.005386739866090491
$ws.Cells.Item(65, 1).Value = ' Make Up and Hair Styling'
$ws.Cells.Item(65, 2).Value = 0.005690079913477174
$ws.Cells.Item(66, 1).Value = 'adaptation'
$ws.Cells.Item(66, 2).Value = 0.005713541297888588
$ws.Cells.Item(67, 1).Value = 'Best Original Screenplay (Spike Jonze)'
$ws.Cells.Item(67, 2).Value = 0.005957026300222758
$ws.Cells.Item(68, 1).Value = '_december'
$ws.Cells.Item(68, 2).Value = 0.006131986705123058
$ws.Cells.Item(69, 1).Value = 'drama'
$ws.Cells.Item(69, 2).Value = 0.006355349585772363
$ws.Cells.Item(70, 1).Value = '_november'
$ws.Cells.Item(70, 2).Value = 0.006532892570951632
$ws.Cells.Item(71, 1).Value = 'Actress in a Leading Role'
$ws.Cells.Item(71, 2).Value = 0.006860597745780897
$ws.Cells.Item(72, 1).Value = ' Best Sound Editing'
$ws.Cells.Item(72, 2).Value = 0.006904829292818856
$ws.Cells.Item(73, 1).Value = ' Film Editing'
$ws.Cells.Item(73, 2).Value = 0.007193430477572829
$ws.Cells.Item(74, 1).Value = ' Ben Wilkins and Thomas Curley)'
$ws.Cells.Item(74, 2).Value = 0.00805735474804927
$ws.Cells.Item(75, 1).Value = 'Best Adapted Screenplay'
$ws.Cells.Item(75, 2).Value = 0.009321696364252462
$ws.Cells.Item(76, 1).Value = 'Best Actor'
$ws.Cells.Item(76, 2).Value = 0.009876108434800527
$ws.Cells.Item(77, 1).Value = 'Best Visual Effects'
$ws.Cells.Item(77, 2).Value = 0.01038297242312005
$ws.Cells.Item(78, 1).Value = 'Unnamed: 0'
$ws.Cells.Item(78, 2).Value = 0.0112231780820014
$ws.Cells.Item(79, 1).Value = 'Best Art Direction'
$ws.Cells.Item(79, 2).Value = 0.01159541573716709
$ws.Cells.Item(80, 1).Value = ' Best Director'
$ws.Cells.Item(80, 2).Value = 0.01165600125880912
$ws.Cells.Item(81, 1).Value = 'history'
$ws.Cells.Item(81, 2).Value = 0.01168590201064186
$ws.Cells.Item(82, 1).Value = ' Best Sound Mixing (Craig Mann'
$ws.Cells.Item(82, 2).Value = 0.01232794899989588
$ws.Cells.Item(83, 1).Value = 'based on a true story'
$ws.Cells.Item(83, 2).Value = 0.01292636317925746
$ws.Cells.Item(84, 1).Value = 'Actress in a Supporting Role'
$ws.Cells.Item(84, 2).Value = 0.01300496570237657
$ws.Cells.Item(85, 1).Value = 'Best Original Song'
$ws.Cells.Item(85, 2).Value = 0.01487460327743167
$ws.Cells.Item(86, 1).Value = 'Best Director'
$ws.Cells.Item(86, 2).Value = 0.01676189614991772
$ws.Cells.Item(87, 1).Value = 'Best Animated Film'
$ws.Cells.Item(87, 2).Value = 0.01766508646193734
$ws.Cells.Item(88, 1).Value = 'Actor in a Supporting Role'
$ws.Cells.Item(88, 2).Value = 0.01943641803447984
$ws.Cells.Item(89, 1).Value = 'average audience'
$ws.Cells.Item(89, 2).Value = 0.01951031273214117
$ws.Cells.Item(90, 1).Value = 'Original Song'
$ws.Cells.Item(90, 2).Value = 0.02060409154272049
$ws.Cells.Item(91, 1).Value = 'Best Original Screenplay'
$ws.Cells.Item(91, 2).Value = 0.02156405796495742
$ws.Cells.Item(92, 1).Value = 'budget recovered'
$ws.Cells.Item(92, 2).Value = 0.02640216761280089
$ws.Cells.Item(93, 1).Value = 'budget recovered opening weekend'
$ws.Cells.Item(93, 2).Value = 0.0272123723899469
$ws.Cells.Item(94, 1).Value = 'imdb rating'
$ws.Cells.Item(94, 2).Value = 0.02834875769072259
$ws.Cells.Item(95, 1).Value = 'Best Supporting Actress'
$ws.Cells.Item(95, 2).Value = 0.03784291968881236
$ws.Cells.Item(96, 1).Value = 'Best Picture'
$ws.Cells.Item(96, 2).Value = 0.03962505381786232
$ws.Cells.Item(97, 1).Value = 'budget ($million)'
$ws.Cells.Item(97, 2).Value = 0.04819187672185874
$ws.Cells.Item(98, 1).Value = 'domestic gross'
$ws.Cells.Item(98, 2).Value = 0.04935033924180245
$ws.Cells.Item(99, 1).Value = 'foreign gross'
$ws.Cells.Item(99, 2).Value = 0.05904429796412305
$ws.Cells.Item(100, 1).Value = 'opening weekend'
$ws.Cells.Item(100, 2).Value = 0.06225877689308793
$ws.Cells.Item(101, 1).Value = 'Best Animated Feature'
$ws.Cells.Item(101, 2).Value = 0.09260779311191843
$ws.Cells.Item(102, 1).Value = 'average critics'
$ws.Cells.Item(102, 2).Value = 0.1584039787904806
